# [FIX] Valeur cohérentes dans la table de paramètrage
#
# The Power Query "bar" (table on sheet "Feuil2", the tab-selected sheet
# backed by worksheets/sheet1.xml) was refreshed from its ODBC source and
# picked up two additional rows. Replay that outcome: append the two new
# rows under the existing query table, then resize the table / named range
# so every dependent part (table ref, autofilter, dimension, defined name)
# stays consistent with the new extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new query rows (IDClient column stays blank, as in the
#     existing rows) -------------------------------------------------------
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 1

$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 0.37999999523162842
$ws.Range("D5").Value = 0.80000001192092896

# --- Grow the "bar" table to cover the refreshed data range --------------
$lo = $ws.ListObjects.Item("bar")
$lo.Resize($ws.Range("A1:D5"))

# --- Re-fit the "hauteur" column now that it holds a longer decimal value -
$ws.Columns.Item(3).AutoFit()

# --- Keep the hidden "DonnéesExternes_1" query-cache name in sync --------
$wb.Names.Item("DonnéesExternes_1").RefersTo = "=Feuil2!`$A`$1:`$D`$5"
